$wb = $excel.ActiveWorkbook

# Update metadata timestamp
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("A2").Value = "05 Nov 2025, 11:39 AM"

# Update "1 Year" column (F) on Industry Analysis sheet
$ia = $wb.Worksheets.Item("Industry Analysis")
$ia.Range("F2").Value = 18.476
$ia.Range("F3").Value = -7.7404
$ia.Range("F4").Value = 30.7972
$ia.Range("F5").Value = -50.2266
$ia.Range("F6").Value = 61.9649
$ia.Range("F7").Value = -9.1713
$ia.Range("F8").Value = -3.556
$ia.Range("F9").Value = 38.3509
$ia.Range("F10").Value = -6.2497
$ia.Range("F11").Value = 52.6723
$ia.Range("F12").Value = -6.932
$ia.Range("F13").Value = 17.5662
$ia.Range("F14").Value = -35.5106
$ia.Range("F15").Value = 0.6286
$ia.Range("F16").Value = -3.1514
$ia.Range("F17").Value = -20.6354
$ia.Range("F18").Value = -0.0175
$ia.Range("F19").Value = -26.9255
$ia.Range("F20").Value = 44.703
$ia.Range("F21").Value = 10.0506
$ia.Range("F22").Value = 84.6016
$ia.Range("F23").Value = -54.4868
$ia.Range("F24").Value = -12.8122
$ia.Range("F25").Value = -9.182700000000001
$ia.Range("F26").Value = 5.9529
$ia.Range("F27").Value = -33.2998
$ia.Range("F28").Value = -20.4441
$ia.Range("F29").Value = -17.1514
$ia.Range("F30").Value = 24.527
$ia.Range("F31").Value = 57.6193
$ia.Range("F32").Value = -1.527
$ia.Range("F33").Value = -5.2378
$ia.Range("F34").Value = 27.4054
$ia.Range("F35").Value = 6.7961
$ia.Range("F36").Value = -5.6683
$ia.Range("F37").Value = 1.4178
$ia.Range("F38").Value = -22.4272
$ia.Range("F39").Value = 12.3741
$ia.Range("F40").Value = -5.138
$ia.Range("F41").Value = -0.1825
$ia.Range("F42").Value = 23.2483
$ia.Range("F43").Value = 14.456
$ia.Range("F44").Value = -11.1739
$ia.Range("F45").Value = 27.112
$ia.Range("F46").Value = -5.6252
$ia.Range("F47").Value = -36.5148
$ia.Range("F48").Value = -27.8397
$ia.Range("F49").Value = -25.4424
$ia.Range("F50").Value = -49.1173
$ia.Range("F51").Value = -51.065
$ia.Range("F52").Value = -35.4517
$ia.Range("F53").Value = -11.9879
$ia.Range("F54").Value = -3.0992
$ia.Range("F55").Value = -15.3441
$ia.Range("F56").Value = -25.937
$ia.Range("F57").Value = -29.1486
$ia.Range("F58").Value = -6.4093
$ia.Range("F59").Value = -23.3046
$ia.Range("F60").Value = -11.2657
$ia.Range("F61").Value = -9.777699999999999
$ia.Range("F62").Value = -16.0561
$ia.Range("F63").Value = -9.932499999999999
$ia.Range("F64").Value = 51.8767
$ia.Range("F65").Value = -43.5191
$ia.Range("F66").Value = 13.7315
$ia.Range("F67").Value = 12.6111
$ia.Range("F68").Value = 31.7532
$ia.Range("F69").Value = -19.9577
$ia.Range("F70").Value = -12.9642
$ia.Range("F71").Value = 13.2432
$ia.Range("F72").Value = 2.8232
$ia.Range("F73").Value = -9.179
$ia.Range("F74").Value = -14.2931
$ia.Range("F75").Value = 28.3699
$ia.Range("F76").Value = 45.5868
